$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 93.32574466666667
$ws.Range("H2").Value = 279.977234
$ws.Range("I2").Value = 0.2327963689879921
$ws.Range("J2").Value = 0.2327963689879922
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 8.131233999999999
$ws.Range("N2").Value = 24.393702
$ws.Range("O2").Value = 0.02090995573015822
$ws.Range("P2").Value = 0.02090995573015823
$ws.Range("Q2").Value = 758.8534681089186
$ws.Range("R2").Value = 6829.681212980267
$ws.Range("S2").Value = 0.004867761769680494
$ws.Range("T2").Value = 0.004867761769680496
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 93.32574466666667
$ws.Range("H3").Value = 279.977234
$ws.Range("I3").Value = 0.2327963689879921
$ws.Range("J3").Value = 0.2327963689879922
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 243.3763986666667
$ws.Range("N3").Value = 730.1291960000001
$ws.Range("O3").Value = 0.625857000534647
$ws.Range("P3").Value = 0.6258570005346471
$ws.Range("Q3").Value = 22713.28363985821
$ws.Range("R3").Value = 204419.5527587239
$ws.Range("S3").Value = 0.1456972372301817
$ws.Range("T3").Value = 0.1456972372301817
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 93.32574466666667
$ws.Range("H4").Value = 279.977234
$ws.Range("I4").Value = 0.2327963689879921
$ws.Range("J4").Value = 0.2327963689879922
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 103.9426383333333
$ws.Range("N4").Value = 311.827915
$ws.Range("O4").Value = 0.2672947262403034
$ws.Range("P4").Value = 0.2672947262403035
$ws.Range("Q4").Value = 9700.524125076345
$ws.Range("R4").Value = 87304.71712568712
$ws.Range("S4").Value = 0.06222524171838201
$ws.Range("T4").Value = 0.06222524171838204
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 93.32574466666667
$ws.Range("H5").Value = 279.977234
$ws.Range("I5").Value = 0.2327963689879921
$ws.Range("J5").Value = 0.2327963689879922
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 33.41874933333333
$ws.Range("N5").Value = 100.256248
$ws.Range("O5").Value = 0.08593831749489127
$ws.Range("P5").Value = 0.08593831749489128
$ws.Range("Q5").Value = 3118.829667362003
$ws.Range("R5").Value = 28069.46700625803
$ws.Range("S5").Value = 0.02000612826974792
$ws.Range("T5").Value = 0.02000612826974793
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 164.7897643333334
$ws.Range("H6").Value = 494.369293
$ws.Range("I6").Value = 0.4110597662007076
$ws.Range("J6").Value = 0.4110597662007077
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 8.131233999999999
$ws.Range("N6").Value = 24.393702
$ws.Range("O6").Value = 0.02090995573015822
$ws.Range("P6").Value = 0.02090995573015823
$ws.Range("Q6").Value = 1339.944134599187
$ws.Range("R6").Value = 12059.49721139269
$ws.Range("S6").Value = 0.008595241513705986
$ws.Range("T6").Value = 0.008595241513705987
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 164.7897643333334
$ws.Range("H7").Value = 494.369293
$ws.Range("I7").Value = 0.4110597662007076
$ws.Range("J7").Value = 0.4110597662007077
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 243.3763986666667
$ws.Range("N7").Value = 730.1291960000001
$ws.Range("O7").Value = 0.625857000534647
$ws.Range("P7").Value = 0.6258570005346471
$ws.Range("Q7").Value = 40105.93938057539
$ws.Range("R7").Value = 360953.4544251785
$ws.Range("S7").Value = 0.2572646323148481
$ws.Range("T7").Value = 0.2572646323148482
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 164.7897643333334
$ws.Range("H8").Value = 494.369293
$ws.Range("I8").Value = 0.4110597662007076
$ws.Range("J8").Value = 0.4110597662007077
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 103.9426383333333
$ws.Range("N8").Value = 311.827915
$ws.Range("O8").Value = 0.2672947262403034
$ws.Range("P8").Value = 0.2672947262403035
$ws.Range("Q8").Value = 17128.6828751349
$ws.Range("R8").Value = 154158.1458762141
$ws.Range("S8").Value = 0.1098741076750213
$ws.Range("T8").Value = 0.1098741076750213
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 164.7897643333334
$ws.Range("H9").Value = 494.369293
$ws.Range("I9").Value = 0.4110597662007076
$ws.Range("J9").Value = 0.4110597662007077
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 33.41874933333333
$ws.Range("N9").Value = 100.256248
$ws.Range("O9").Value = 0.08593831749489127
$ws.Range("P9").Value = 0.08593831749489128
$ws.Range("Q9").Value = 5507.067826954741
$ws.Range("R9").Value = 49563.61044259267
$ws.Range("S9").Value = 0.03532578469713218
$ws.Range("T9").Value = 0.0353257846971322
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 95.835818
$ws.Range("H10").Value = 287.507454
$ws.Range("I10").Value = 0.2390576204784642
$ws.Range("J10").Value = 0.2390576204784643
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.131233999999999
$ws.Range("N10").Value = 24.393702
$ws.Range("O10").Value = 0.02090995573015822
$ws.Range("P10").Value = 0.02090995573015823
$ws.Range("Q10").Value = 779.2634617394119
$ws.Range("R10").Value = 7013.371155654707
$ws.Range("S10").Value = 0.004998684261161653
$ws.Range("T10").Value = 0.004998684261161655
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 95.835818
$ws.Range("H11").Value = 287.507454
$ws.Range("I11").Value = 0.2390576204784642
$ws.Range("J11").Value = 0.2390576204784643
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 243.3763986666667
$ws.Range("N11").Value = 730.1291960000001
$ws.Range("O11").Value = 0.625857000534647
$ws.Range("P11").Value = 0.6258570005346471
$ws.Range("Q11").Value = 23324.17624811411
$ws.Range("R11").Value = 209917.586233027
$ws.Range("S11").Value = 0.1496158853076016
$ws.Range("T11").Value = 0.1496158853076017
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 95.835818
$ws.Range("H12").Value = 287.507454
$ws.Range("I12").Value = 0.2390576204784642
$ws.Range("J12").Value = 0.2390576204784643
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 103.9426383333333
$ws.Range("N12").Value = 311.827915
$ws.Range("O12").Value = 0.2672947262403034
$ws.Range("P12").Value = 0.2672947262403035
$ws.Range("Q12").Value = 9961.427769753158
$ws.Range("R12").Value = 89652.84992777841
$ws.Range("S12").Value = 0.06389884122144945
$ws.Range("T12").Value = 0.06389884122144947
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 95.835818
$ws.Range("H13").Value = 287.507454
$ws.Range("I13").Value = 0.2390576204784642
$ws.Range("J13").Value = 0.2390576204784643
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 33.41874933333333
$ws.Range("N13").Value = 100.256248
$ws.Range("O13").Value = 0.08593831749489127
$ws.Range("P13").Value = 0.08593831749489128
$ws.Range("Q13").Value = 3202.713178896955
$ws.Range("R13").Value = 28824.41861007259
$ws.Range("S13").Value = 0.02054420968825148
$ws.Range("T13").Value = 0.02054420968825148
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 46.93870866666666
$ws.Range("H14").Value = 140.816126
$ws.Range("I14").Value = 0.117086244332836
$ws.Range("J14").Value = 0.117086244332836
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 8.131233999999999
$ws.Range("N14").Value = 24.393702
$ws.Range("O14").Value = 0.02090995573015822
$ws.Range("P14").Value = 0.02090995573015823
$ws.Range("Q14").Value = 381.6696238264946
$ws.Range("R14").Value = 3435.026614438451
$ws.Range("S14").Value = 0.002448268185610089
$ws.Range("T14").Value = 0.00244826818561009
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 46.93870866666666
$ws.Range("H15").Value = 140.816126
$ws.Range("I15").Value = 0.117086244332836
$ws.Range("J15").Value = 0.117086244332836
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 243.3763986666667
$ws.Range("N15").Value = 730.1291960000001
$ws.Range("O15").Value = 0.625857000534647
$ws.Range("P15").Value = 0.6258570005346471
$ws.Range("Q15").Value = 11423.77387335719
$ws.Range("R15").Value = 102813.9648602147
$ws.Range("S15").Value = 0.07327924568201552
$ws.Range("T15").Value = 0.07327924568201555
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 46.93870866666666
$ws.Range("H16").Value = 140.816126
$ws.Range("I16").Value = 0.117086244332836
$ws.Range("J16").Value = 0.117086244332836
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 103.9426383333333
$ws.Range("N16").Value = 311.827915
$ws.Range("O16").Value = 0.2672947262403034
$ws.Range("P16").Value = 0.2672947262403035
$ws.Range("Q16").Value = 4878.933218773032
$ws.Range("R16").Value = 43910.39896895729
$ws.Range("S16").Value = 0.03129653562545066
$ws.Range("T16").Value = 0.03129653562545068
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 46.93870866666666
$ws.Range("H17").Value = 140.816126
$ws.Range("I17").Value = 0.117086244332836
$ws.Range("J17").Value = 0.117086244332836
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 33.41874933333333
$ws.Range("N17").Value = 100.256248
$ws.Range("O17").Value = 0.08593831749489127
$ws.Range("P17").Value = 0.08593831749489128
$ws.Range("Q17").Value = 1568.632938961694
$ws.Range("R17").Value = 14117.69645065525
$ws.Range("S17").Value = 0.01006219483975967
$ws.Range("T17").Value = 0.01006219483975967